$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update User_ID column (A) and Rating column (C) per the new data values
$ws.Range("A2").Value = 8013
$ws.Range("C2").Value = 2

$ws.Range("A3").Value = 8941
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 1581
$ws.Range("C4").Value = 4

$ws.Range("A5").Value = 5719
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = 7040
